$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "primary_procedure" column (L) is being removed from the dataset;
# deleting it shifts every subsequent column one place to the left.
$ws.Columns("L").Delete()

# Re-apply the column widths the author set (on the post-delete layout)
# so the remaining columns display their contents comfortably.
$ws.Columns("E").ColumnWidth = 21.998697916666668
$ws.Columns("F").ColumnWidth = 16.666666666666668
$ws.Columns("G").ColumnWidth = 16.998697916666668
$ws.Columns("J").ColumnWidth = 15.166666666666668
$ws.Columns("K").ColumnWidth = 18.998697916666668
$ws.Columns("L").ColumnWidth = 13.830729166666668
$ws.Columns("N").ColumnWidth = 13.166666666666668
$ws.Columns("P").ColumnWidth = 13.998697916666668

# Leave the active selection where the author left it after the edit.
$ws.Range("K13").Select() | Out-Null
